# Updated Bug Tracker With Fix
# Adds a "Date Fixed" and "Fix Description" entry for bug #1 on the "Bugs" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bugs")

# Date Fixed (column G) - use the same date format already applied to the
# "Date Found" column (E), which is the custom "dd/mm/yyyy;@" format.
$ws.Range("G3").Value2 = 45528
$ws.Range("G3").NumberFormat = $ws.Range("E3").NumberFormat

# Fix Description (column H)
$ws.Range("H3").Value2 = "Use server travel on session destroyed delegate"

# Widen column H so the new fix description is fully visible.
$ws.Columns.Item(8).ColumnWidth = 40

# Leave the new fix-description cell selected, as in the saved workbook.
$ws.Range("H3").Select() | Out-Null
